$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.228422045707703
$ws.Range("B1").Value = 1.98248028755188
$ws.Range("C1").Value = 4.308877468109131
$ws.Range("D1").Value = 2.987842559814453
$ws.Range("E1").Value = 1.181873798370361
